{"js": "// Replace the text of each arithmetic-problem cell in the single table,\n// in row-major document order, with its new value \u2014 preserving each\n// cell's existing paragraph/run formatting (font, size, alignment).\nconst newValues = [\"42-39=\", \"85-76=\", \"49+32=\", \"90-7=\", \"46+38=\", \"29+66=\", \"19+15=\", \"69+3=\", \"40-9=\", \"76+6=\", \"39+57=\", \"8+88=\", \"84-68=\", \"93-88=\", \"15+19=\", \"70-9=\", \"47+47=\", \"73-18=\", \"48+38=\", \"97-49=\", \"57-8=\", \"54+27=\", \"29+24=\", \"75-16=\", \"58+13=\", \"62-56=\", \"50-29=\", \"76-39=\", \"58+37=\", \"33-9=\", \"56+9=\", \"33-27=\", \"73-68=\", \"73-17=\", \"62-34=\", \"30-22=\", \"86-27=\", \"82-43=\", \"57-9=\", \"29+16=\", \"5+78=\", \"87+4=\", \"28+46=\", \"64-48=\", \"92-68=\", \"56-48=\", \"92-33=\", \"73-34=\", \"19+38=\", \"89+8=\", \"93-7=\", \"4+19=\", \"16+45=\", \"80-78=\", \"64+27=\", \"91-22=\", \"23+8=\", \"72-17=\", \"51-23=\", \"91-66=\", \"59+37=\", \"53+8=\", \"7+67=\", \"52-16=\", \"72-64=\", \"8+6=\", \"6+28=\", \"83-29=\", \"3+48=\", \"82-43=\", \"22-4=\", \"19+17=\", \"83-58=\", \"44-6=\", \"26+57=\", \"92-35=\", \"67+14=\", \"72-15=\", \"7+38=\", \"71-49=\", \"34+29=\", \"92-58=\", \"19+48=\", \"47+6=\", \"6+9=\", \"18+48=\", \"90-53=\", \"13-9=\", \"36+57=\", \"65+17=\", \"29+63=\", \"58+28=\", \"55+8=\", \"8+14=\", \"64-57=\", \"44-19=\", \"70-9=\", \"41-37=\", \"95-38=\", \"30-21=\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Determine column count from the first row's cells.\nconst firstRowCells = rows.items[0].cells;\nfirstRowCells.load(\"items\");\nawait context.sync();\nconst colCount = firstRowCells.items.length;\n\nlet idx = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    if (idx >= newValues.length) break;\n    const cell = table.getCell(r, c);\n    const paras = cell.body.paragraphs;\n    paras.load(\"items\");\n    await context.sync();\n    const p = paras.items[0];\n    const range = p.getRange();\n    range.insertText(newValues[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the text of each arithmetic-problem cell in the single table,\n# in row-major document order, with its new value -- preserving each\n# cell's existing paragraph/run formatting (font, size, alignment).\n$newValues = @(\"42-39=\", \"85-76=\", \"49+32=\", \"90-7=\", \"46+38=\", \"29+66=\", \"19+15=\", \"69+3=\", \"40-9=\", \"76+6=\", \"39+57=\", \"8+88=\", \"84-68=\", \"93-88=\", \"15+19=\", \"70-9=\", \"47+47=\", \"73-18=\", \"48+38=\", \"97-49=\", \"57-8=\", \"54+27=\", \"29+24=\", \"75-16=\", \"58+13=\", \"62-56=\", \"50-29=\", \"76-39=\", \"58+37=\", \"33-9=\", \"56+9=\", \"33-27=\", \"73-68=\", \"73-17=\", \"62-34=\", \"30-22=\", \"86-27=\", \"82-43=\", \"57-9=\", \"29+16=\", \"5+78=\", \"87+4=\", \"28+46=\", \"64-48=\", \"92-68=\", \"56-48=\", \"92-33=\", \"73-34=\", \"19+38=\", \"89+8=\", \"93-7=\", \"4+19=\", \"16+45=\", \"80-78=\", \"64+27=\", \"91-22=\", \"23+8=\", \"72-17=\", \"51-23=\", \"91-66=\", \"59+37=\", \"53+8=\", \"7+67=\", \"52-16=\", \"72-64=\", \"8+6=\", \"6+28=\", \"83-29=\", \"3+48=\", \"82-43=\", \"22-4=\", \"19+17=\", \"83-58=\", \"44-6=\", \"26+57=\", \"92-35=\", \"67+14=\", \"72-15=\", \"7+38=\", \"71-49=\", \"34+29=\", \"92-58=\", \"19+48=\", \"47+6=\", \"6+9=\", \"18+48=\", \"90-53=\", \"13-9=\", \"36+57=\", \"65+17=\", \"29+63=\", \"58+28=\", \"55+8=\", \"8+14=\", \"64-57=\", \"44-19=\", \"70-9=\", \"41-37=\", \"95-38=\", \"30-21=\")\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    if ($idx -ge $newValues.Length) { break }\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$idx]\n    $idx++\n  }\n}\n"}
